$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string (label) edits ---
# A8: "Volume 30   Number  13" -> "...14" (issue number)
$ws.Range("A8").Characters(21,2).Text = "14"

# C9: "Report Covering the Week  3/27/2023  Through  4/2/2023" -> updated week dates
$ws.Range("C9").Characters(27,9).Text = "4/3/2023"
$ws.Range("C9").Characters(46,8).Text = "4/9/2023"

# --- Data table edits (rows 15-30) ---

$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("E15").Value = 0
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -33.333333333333
$ws.Range("I15").Value = 7
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = -12.5
$ws.Range("L15").Value = -22.222222222222
$ws.Range("M15").Value = 40
$ws.Range("N15").Value = -12.5
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 17.647058823529
$ws.Range("I16").Value = 66
$ws.Range("J16").Value = 67
$ws.Range("K16").Value = -1.492537313432
$ws.Range("L16").Value = 17.857142857142
$ws.Range("M16").Value = -24.137931034482
$ws.Range("N16").Value = -80.924855491329
$ws.Range("C17").Value = 5
$ws.Range("E17").Value = -37.5
$ws.Range("F17").Value = 34
$ws.Range("G17").Value = 36
$ws.Range("H17").Value = -5.555555555555
$ws.Range("I17").Value = 124
$ws.Range("J17").Value = 109
$ws.Range("K17").Value = 13.761467889908
$ws.Range("L17").Value = 53.086419753086
$ws.Range("M17").Value = 58.974358974359
$ws.Range("N17").Value = 44.186046511627
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -23.076923076923
$ws.Range("I18").Value = 30
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = -28.571428571428
$ws.Range("L18").Value = -23.076923076923
$ws.Range("M18").Value = -63.855421686747
$ws.Range("N18").Value = -95.398773006135
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -31.25
$ws.Range("F19").Value = 56
$ws.Range("G19").Value = 61
$ws.Range("H19").Value = -8.196721311475
$ws.Range("I19").Value = 208
$ws.Range("J19").Value = 304
$ws.Range("K19").Value = -31.578947368421
$ws.Range("L19").Value = 52.941176470588
$ws.Range("M19").Value = 80.869565217391
$ws.Range("N19").Value = -36.969696969697
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 20
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = -16.666666666666
$ws.Range("I20").Value = 80
$ws.Range("J20").Value = 69
$ws.Range("K20").Value = 15.942028985507
$ws.Range("L20").Value = 100
$ws.Range("M20").Value = 29.032258064516
$ws.Range("N20").Value = -86.577181208053
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = -36.585365853658
$ws.Range("F21").Value = 142
$ws.Range("G21").Value = 154
$ws.Range("H21").Value = -7.792207792207
$ws.Range("I21").Value = 515
$ws.Range("J21").Value = 601
$ws.Range("K21").Value = -14.309484193011
$ws.Range("L21").Value = 42.265193370165
$ws.Range("M21").Value = 19.489559164733
$ws.Range("N21").Value = -74.567901234567
$ws.Range("C22").Value = 2
$ws.Range("I22").Value = 30
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = 172.727272727273
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = 3.571428571428
$ws.Range("F24").Value = 137
$ws.Range("G24").Value = 105
$ws.Range("H24").Value = 30.476190476190
$ws.Range("I24").Value = 583
$ws.Range("J24").Value = 441
$ws.Range("K24").Value = 32.199546485260
$ws.Range("L24").Value = 79.938271604938
$ws.Range("M24").Value = 132.270916334661
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = 5.882352941176
$ws.Range("F25").Value = 76
$ws.Range("G25").Value = 49
$ws.Range("H25").Value = 55.102040816326
$ws.Range("I25").Value = 235
$ws.Range("J25").Value = 238
$ws.Range("K25").Value = -1.260504201680
$ws.Range("L25").Value = 35.838150289017
$ws.Range("M25").Value = -7.114624505928
$ws.Range("C26").Value = 2
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 5
$ws.Range("H26").Value = 66.666666666666
$ws.Range("I26").Value = 15
$ws.Range("J26").Value = 12
$ws.Range("K26").Value = 25
$ws.Range("L26").Value = 15.384615384615
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -33.333333333333
$ws.Range("F27").Value = 10
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 66.666666666666
$ws.Range("I27").Value = 29
$ws.Range("J27").Value = 25
$ws.Range("K27").Value = 16
$ws.Range("L27").Value = -9.375
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = '#,##0'
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 1
$ws.Range("G30").NumberFormat = '#,##0'
$ws.Range("H30").Value = 0
$ws.Range("H30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J30").Value = 2
$ws.Range("K30").Value = 50
